$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.038511753082275
$ws.Range("B1").Value = 3.474916696548462
$ws.Range("C1").Value = 2.690839290618896
$ws.Range("D1").Value = 2.569641351699829
$ws.Range("E1").Value = 2.614851236343384
